$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.522.48'
$ws.Range("E2").Value = '  +3.67%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.816.82'
$ws.Range("E3").Value = '  +4.39%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '343.51'
$ws.Range("E5").Value = '  +3.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9991'
$ws.Range("E6").Value = '  -0.18%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3825'
$ws.Range("E7").Value = '  +2.22%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3556'
$ws.Range("E8").Value = '  +4.32%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '50.40'
$ws.Range("E9").Value = '  +1.85%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.237'
$ws.Range("E10").Value = '  +3.45%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07789'
$ws.Range("E11").Value = '  +4.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.001'
$ws.Range("E12").Value = '  -0.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.41'
$ws.Range("E13").Value = '  +10.22%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.606'
$ws.Range("E14").Value = '  +2.43%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.814.37'
$ws.Range("E15").Value = '  +4.10%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.215'
$ws.Range("E16").Value = '  +1.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001131'
$ws.Range("E17").Value = '  +4.03%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06724'
$ws.Range("E18").Value = '  +0.54%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '87.13'
$ws.Range("E19").Value = '  +4.11%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("E20").Value = '  -0.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.64'
$ws.Range("E21").Value = '  +5.80%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.552'
$ws.Range("E22").Value = '  +5.96%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.22'
$ws.Range("E23").Value = '  +1.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '27.536.69'
$ws.Range("E24").Value = '  +3.51%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.478'
$ws.Range("E25").Value = '  +0.60%  '

$ws.Range("E26").Value = '  +9.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.22'
$ws.Range("E27").Value = '  +13.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.464'
$ws.Range("E28").Value = '  +4.78%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '153.52'
$ws.Range("E29").Value = '  +0.58%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.018.50'
$ws.Range("E30").Value = '  +3.85%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '135.74'
$ws.Range("E31").Value = '  +2.62%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.393'
$ws.Range("E32").Value = '  +4.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.100'
$ws.Range("E33").Value = '  -0.59%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.89'
$ws.Range("E34").Value = '  +6.82%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08767'
$ws.Range("E35").Value = '  +2.07%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.699'
$ws.Range("E36").Value = '  -1.21%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.660'
$ws.Range("E37").Value = '  +3.95%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.7059'
$ws.Range("E38").Value = '  +13.19%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.108'
$ws.Range("E39").Value = '  +5.73%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06521'
$ws.Range("E40").Value = '  +3.39%  '

$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2257'
$ws.Range("E41").Value = '  +3.87%  '

$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.02410'
$ws.Range("E42").Value = '  +2.13%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.316'
$ws.Range("E43").Value = '  +6.70%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.80'
$ws.Range("E44").Value = '  +3.02%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6682'
$ws.Range("E45").Value = '  +10.45%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9992'
$ws.Range("E46").Value = '  -0.31%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.952'
$ws.Range("E47").Value = '  +1.20%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.192'
$ws.Range("E48").Value = '  +6.22%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '132.68'
$ws.Range("E49").Value = '  +2.78%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07312'
$ws.Range("E50").Value = '  +0.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.73'
$ws.Range("E51").Value = '  +3.78%  '
